$wb = $excel.ActiveWorkbook

# --- Update the conversion text on sheet "Hoja1" (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.04 = 7564.42 pesos`n✅ 7564.42 pesos = 2.03 = 922.07 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update numeric rate cells on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("O10").Value = 3699
$wsTasas.Range("N12").Value = 3732.69
$wsTasas.Range("O12").Value = 455
